$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: MIT Governance Lab (GOV/LAB) now matched with AIR-INK
$ws.Range("B27").Value = "['None'],AIR-INK: Air-Pollution to ink"
$ws.Range("C27").Value = 1

# Row 31: Putnam Associates no longer matched with AIR-INK
$ws.Range("B31").Value = "['None']"
$ws.Range("C31").Value = 0

# Row 38: The Kroger Co. Zero Hunger Zero Waste Foundation matched with Algramo
$ws.Range("B38").Value = "['None'],Algramo-Catalyzing Reusable Packaging"
$ws.Range("C38").Value = 1
